# Insert a new weekly price report (2 rows) for
# "Feria Lagunitas de Puerto Montt" / Plátano, just above the row that used
# to be row 663, shifting the remainder of the table down by two rows.
#
# Net effect matches the target diff:
#   - dimension grows from A1:T766 to A1:T768
#   - the two newly inserted rows (663 & 664) carry the same
#     Mercado/Region/Categoria/Variedad/Calidad/Volumen/Unidad/Origen/Kg-unidad
#     values as the (old) rows 663 & 664, but with an updated Fecha
#     (44984) and updated price columns (N, O, P, S).
#   - everything that used to occupy rows 663-766 simply shifts down to
#     665-768 untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 663:664 - this shifts the existing rows 663-766
# down to 665-768 (and carries the row-663/664 formatting, e.g. the date
# number-format on column D, onto the freshly inserted rows).
$ws.Range("663:664").Insert()

# --- New row 663 (copy of old row 663's static columns, new Fecha + prices) ---
$ws.Cells.Item(663,1).Value2  = 4
$ws.Cells.Item(663,2).Value2  = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(663,3).Value2  = 'Los Lagos'
$ws.Cells.Item(663,4).Value2  = 44984
$ws.Cells.Item(663,5).Value2  = 10
$ws.Cells.Item(663,6).Value2  = 'Fruta'
$ws.Cells.Item(663,7).Value2  = 100108
$ws.Cells.Item(663,8).Value2  = 'Tropicales y subtropicales'
$ws.Cells.Item(663,9).Value2  = 100108006
$ws.Cells.Item(663,10).Value2 = 'Plátano'
$ws.Cells.Item(663,11).Value2 = 'Sin especificar'
$ws.Cells.Item(663,12).Value2 = 'Pintón'
$ws.Cells.Item(663,13).Value2 = 400
$ws.Cells.Item(663,14).Value2 = 26000
$ws.Cells.Item(663,15).Value2 = 26000
$ws.Cells.Item(663,16).Value2 = 26000
$ws.Cells.Item(663,17).Value2 = '$/caja 20 kilos'
$ws.Cells.Item(663,18).Value2 = 'Ecuador'
$ws.Cells.Item(663,19).Value2 = 1300
$ws.Cells.Item(663,20).Value2 = 20

# --- New row 664 (copy of old row 664's static columns, new Fecha + prices) ---
$ws.Cells.Item(664,1).Value2  = 4
$ws.Cells.Item(664,2).Value2  = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(664,3).Value2  = 'Los Lagos'
$ws.Cells.Item(664,4).Value2  = 44984
$ws.Cells.Item(664,5).Value2  = 10
$ws.Cells.Item(664,6).Value2  = 'Fruta'
$ws.Cells.Item(664,7).Value2  = 100108
$ws.Cells.Item(664,8).Value2  = 'Tropicales y subtropicales'
$ws.Cells.Item(664,9).Value2  = 100108006
$ws.Cells.Item(664,10).Value2 = 'Plátano'
$ws.Cells.Item(664,11).Value2 = 'Sin especificar'
$ws.Cells.Item(664,12).Value2 = 'Primera Pintón'
$ws.Cells.Item(664,13).Value2 = 800
$ws.Cells.Item(664,14).Value2 = 28000
$ws.Cells.Item(664,15).Value2 = 29000
$ws.Cells.Item(664,16).Value2 = 28500
$ws.Cells.Item(664,17).Value2 = '$/caja 20 kilos'
$ws.Cells.Item(664,18).Value2 = 'Ecuador'
$ws.Cells.Item(664,19).Value2 = 1425
$ws.Cells.Item(664,20).Value2 = 20
